# Generate Report for Handoff
# - Updates the localization-status report so the zh-cn / de-de rows (and the
#   Overview roll-up) reflect that the content is now "Ready for handoff"
#   instead of "Handed back: in sync with en-US", and refreshes the
#   handoff/generation timestamps accordingly.
# - The Status column is much narrower now that the text is shorter, so the
#   column width is refreshed to fit the new text (mirrors an Excel
#   AutoFit-style resize after the text shrank).

$wb = $excel.ActiveWorkbook

$statusText = "Ready for handoff"

# Column width (character units) that this host's rounding produces the
# closest match to the real workbook's post-resize width of ~17.216.
$newStatusColWidth = 16.333333333333332

# ---------------------------------------------------------------------
# Overview sheet: zh-cn / de-de status + latest handoff-generation date
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $statusText
$wsOverview.Range("F2").Value = $statusText
$wsOverview.Range("G2").Value = "2016-09-06 15:43:36"

$wsOverview.Columns.Item(5).ColumnWidth = $newStatusColWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newStatusColWidth

# ---------------------------------------------------------------------
# zh-cn sheet: Status + Latest Handoff Datetime
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $statusText
$wsZhCn.Range("H2").Value = "2016-09-06 15:43:23"

$wsZhCn.Columns.Item(3).ColumnWidth = $newStatusColWidth

# ---------------------------------------------------------------------
# de-de sheet: Status + Latest Handoff Datetime
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $statusText
$wsDeDe.Range("H2").Value = "2016-09-06 15:43:36"

$wsDeDe.Columns.Item(3).ColumnWidth = $newStatusColWidth
